# Automatische test-sync: 2025-08-05 17:49:50
# Append the new "Kun jij dit even regelen?" mail-log row (row 19) to the
# Logs sheet, extend the conditional-formatting ranges that track the
# sheet's used range, and bump the Dashboard's "Planning / Afspraak" count.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new log row -------------------------------------------
$row = 19
$logs.Cells.Item($row, 1).Value = "Kun jij dit even regelen?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Cells.Item($row, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($row, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($row, 6).Value = "2025-08-05 17:49:15"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"
$logs.Cells.Item($row, 10).Value = "Nee"

# --- 2. Extend the conditional formatting ranges to include row 19 -------
$logs.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D19"))
$logs.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G19"))
$logs.Range("H2:H18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H19"))
$logs.Range("I2:I18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I19"))
$logs.Range("J2:J18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J19"))

# --- 3. Bump the Dashboard count for "Planning / Afspraak" ---------------
$dashboard.Range("B2").Value = 13
